$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D and add a new column E with data for each person's
# second week of dedication hours.
$ws.Columns.Item(4).ColumnWidth = 96.83333333333333
$ws.Columns.Item(5).ColumnWidth = 55.5

# Row 2 - Ivan
$ws.Range("D2").Value = "1'5h(organización reuniones, documentacion)"
$ws.Range("E2").Value = "2'5h(reuniones, actas y documentación)"

# Row 3 - Marta
$ws.Range("E3").Value = "2h(maquetación y documentación)"

# Row 4 - Alfonso
$ws.Range("D4").Value = "6h(reunión,implementación web y BBDD)"
$ws.Range("E4").Value = "4h(Implementación, poblado de BBDD y pruebas)"

# Row 5 - Sandra
$ws.Range("D5").Value = "2h(documentacion)"
$ws.Range("E5").Value = "2h(reuniones, documentacion)"

# Row 6 - Dariel
$ws.Range("D6").Value = "8h(reunion, diagramas, readme, modulo html detalles, método get de detalles, funcionalidad de busquedas)"
$ws.Range("E6").Value = "6h(reunion,funcionalidades de análisis, comentarios, bug fix y pruebas)"

# Row 7 - Leon
$ws.Range("E7").Value = "3h(reunion y documentacion)"

# Row 8 - Alberto
$ws.Range("D8").Value = "1'5h(organización reuniones, documentacion)"
$ws.Range("E8").Value = "1'5h(organización reuniones, documentacion)"

# Row 9 - M.Angel
$ws.Range("D9").Value = "1h(Documentación)"
$ws.Range("E9").Value = "2h(reunion y documentacion)"

# Row 10 - Bea
$ws.Range("E10").Value = "2'5h(reuniones, actas y documentación)"

# Update the selection to match where the author ended up (E7) after
# entering the new week's data.
$ws.Range("E7").Select()

$wb.Save()
